# Insert a new data row right above the current row 37 (Fecha=2022-12-28 /
# serial 44923), shifting all the rows below it down by one, then fill the
# newly inserted row with its own values (row 37 becomes a brand-new record;
# everything that used to be row 37 onward just moves down to row 38 onward).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 - this shifts rows 37:110 down to 38:111
# and copies formatting (e.g. the date style in column D) from the row above.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 45002
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 100112022
$ws.Cells.Item(37, 7).Value = "Arveja Verde"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 20
$ws.Cells.Item(37, 11).Value = 35000
$ws.Cells.Item(37, 12).Value = 35000
$ws.Cells.Item(37, 13).Value = 35000
$ws.Cells.Item(37, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(37, 16).Value = 1400
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"
